$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 422.83334
$ws.Range("I5").Value = 336.25
$ws.Range("K5").Value = 336.25
$ws.Range("M5").Value = -221.25

$ws.Range("H9").Value = 716.2
$ws.Range("I9").Value = 25.666666
$ws.Range("J9").Value = 1752
$ws.Range("K9").Value = 25.666666
$ws.Range("L9").Value = 1752
$ws.Range("M9").Value = 143.333334
$ws.Range("N9").Value = -2090

$ws.Range("H39").Value = 283
$ws.Range("I39").Value = 85
$ws.Range("K39").Value = 255
$ws.Range("M39").Value = 41

$ws.Range("H40").Value = 5871.4287
$ws.Range("J40").Value = 6566.6665
$ws.Range("L40").Value = 6566.6665
$ws.Range("N40").Value = -6916.6665

$ws.Range("H70").Value = 1857.5
$ws.Range("J70").Value = 3047.818
$ws.Range("L70").Value = 9143.454000000002
$ws.Range("N70").Value = -9683.454000000002

$ws.Range("H73").Value = 1857.5
$ws.Range("J73").Value = 3047.818
$ws.Range("L73").Value = 9143.454000000002
$ws.Range("N73").Value = -11015.454

$ws.Range("H80").Value = 4425.8335
$ws.Range("I80").Value = 5762.625
$ws.Range("K80").Value = 17287.875
$ws.Range("M80").Value = -16289.875

$ws.Range("H83").Value = 4425.8335
$ws.Range("I83").Value = 5762.625
$ws.Range("K83").Value = 51863.625
$ws.Range("M83").Value = -46871.625

$ws.Range("H113").Value = 4034.2727
$ws.Range("I113").Value = 4381.8887
$ws.Range("J113").Value = 2470
$ws.Range("K113").Value = 4381.8887
$ws.Range("L113").Value = 2470
$ws.Range("M113").Value = -1127.8887
$ws.Range("N113").Value = -8978

$ws.Range("H137").Value = 5046.909
$ws.Range("I137").Value = 4158.421
$ws.Range("K137").Value = 12475.263
$ws.Range("M137").Value = -9925.263000000001

$ws.Range("H138").Value = 5511.8486
$ws.Range("J138").Value = 5895.396
$ws.Range("L138").Value = 17686.188
$ws.Range("N138").Value = -27966.188

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12856.95
$ws.Range("I32").Value = 9180.882
$ws.Range("J32").Value = 24497.834
$ws.Range("K32").Value = 9180.882
$ws.Range("L32").Value = 24497.834
$ws.Range("M32").Value = -8893.882
$ws.Range("N32").Value = -25071.834

$ws.Range("H74").Value = 3985.6316
$ws.Range("I74").Value = 3761.1875
$ws.Range("K74").Value = 3761.1875
$ws.Range("M74").Value = -2887.1875

$ws.Range("H77").Value = 3985.6316
$ws.Range("I77").Value = 3761.1875
$ws.Range("K77").Value = 18805.9375
$ws.Range("M77").Value = -14437.9375

$ws.Range("H139").Value = 74474.625
$ws.Range("J139").Value = 74474.625
$ws.Range("L139").Value = 74474.625
$ws.Range("N139").Value = -84754.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2611.5938
$ws.Range("I20").Value = 2498.9473
$ws.Range("K20").Value = 2498.9473
$ws.Range("M20").Value = -2251.9473

$ws.Range("H86").Value = 1888.4117
$ws.Range("I86").Value = 2168.7273
$ws.Range("K86").Value = 2168.7273
$ws.Range("M86").Value = -1045.7273

$ws.Range("H89").Value = 1888.4117
$ws.Range("I89").Value = 2168.7273
$ws.Range("K89").Value = 10843.6365
$ws.Range("M89").Value = -5227.636500000001

$ws.Range("H135").Value = 78789.81
$ws.Range("J135").Value = 77780.625
$ws.Range("L135").Value = 77780.625
$ws.Range("N135").Value = -87920.625

$ws.Range("H138").Value = 94434.91
$ws.Range("J138").Value = 94434.91
$ws.Range("L138").Value = 94434.91
$ws.Range("N138").Value = -104714.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 666.55554
$ws.Range("I22").Value = 372.85715
$ws.Range("K22").Value = 372.85715
$ws.Range("M22").Value = -22.85714999999999

$ws.Range("H31").Value = 4185.806
$ws.Range("I31").Value = 3010.5334
$ws.Range("J31").Value = 4524.827
$ws.Range("K31").Value = 3010.5334
$ws.Range("L31").Value = 4524.827
$ws.Range("M31").Value = -2715.5334
$ws.Range("N31").Value = -5114.827

$ws.Range("H34").Value = 4185.806
$ws.Range("I34").Value = 3010.5334
$ws.Range("J34").Value = 4524.827
$ws.Range("K34").Value = 3010.5334
$ws.Range("L34").Value = 4524.827
$ws.Range("M34").Value = -2808.5334
$ws.Range("N34").Value = -4928.827

$ws.Range("H51").Value = 41713.855
$ws.Range("J51").Value = 41713.855
$ws.Range("L51").Value = 41713.855
$ws.Range("N51").Value = -43185.855

$ws.Range("H61").Value = 41713.855
$ws.Range("J61").Value = 41713.855
$ws.Range("L61").Value = 41713.855
$ws.Range("N61").Value = -42409.855

$ws.Range("H62").Value = 10207.786
$ws.Range("I62").Value = 12180.8
$ws.Range("K62").Value = 12180.8
$ws.Range("M62").Value = -11556.8

$ws.Range("H65").Value = 10207.786
$ws.Range("I65").Value = 12180.8
$ws.Range("K65").Value = 60904
$ws.Range("M65").Value = -57784

$ws.Range("H69").Value = 12091
$ws.Range("I69").Value = 12091
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 12091
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -11342
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 12091
$ws.Range("I72").Value = 12091
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 36273
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -32529
$ws.Range("N72").ClearContents()

$ws.Range("H141").Value = 355925.97
$ws.Range("J141").Value = 402810.44
$ws.Range("L141").Value = 402810.44
$ws.Range("N141").Value = -413170.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 435.13333
$ws.Range("I2").Value = 88.25
$ws.Range("J2").Value = 831.5714
$ws.Range("K2").Value = 529.5
$ws.Range("L2").Value = 4989.428400000001
$ws.Range("M2").Value = -416.5
$ws.Range("N2").Value = -5215.428400000001

$ws.Range("H107").Value = 1463.7587
$ws.Range("I107").Value = 206.33333
$ws.Range("K107").Value = 618.99999
$ws.Range("M107").Value = 1301.00001

$ws.Range("H133").Value = 13295
$ws.Range("I133").Value = 8286.666999999999
$ws.Range("J133").Value = 16300
$ws.Range("K133").Value = 24860.001
$ws.Range("L133").Value = 48900
$ws.Range("M133").Value = -19800.001
$ws.Range("N133").Value = -59020

$ws.Range("H134").Value = 11281.75
$ws.Range("I134").Value = 15251.8
$ws.Range("J134").Value = 9477.182000000001
$ws.Range("K134").Value = 45755.39999999999
$ws.Range("L134").Value = 28431.546
$ws.Range("M134").Value = -40685.39999999999
$ws.Range("N134").Value = -38571.546

$ws.Range("H136").Value = 3691.4119
$ws.Range("I136").Value = 3483.5715
$ws.Range("K136").Value = 10450.7145
$ws.Range("M136").Value = -5350.7145

$ws.Range("H138").Value = 7883.8335
$ws.Range("J138").Value = 8992.444
$ws.Range("L138").Value = 26977.332
$ws.Range("N138").Value = -37257.33199999999

$ws.Range("H139").Value = 6349
$ws.Range("I139").Value = 7315.4287
$ws.Range("J139").Value = 2966.5
$ws.Range("K139").Value = 21946.2861
$ws.Range("L139").Value = 8899.5
$ws.Range("M139").Value = -16806.2861
$ws.Range("N139").Value = -19179.5

$ws.Range("H140").Value = 3126.125
$ws.Range("I140").Value = 2001.8
$ws.Range("K140").Value = 6005.4
$ws.Range("M140").Value = -825.3999999999996

$ws.Range("H141").Value = 6234.2
$ws.Range("I141").Value = 6234.2
$ws.Range("K141").Value = 18702.6
$ws.Range("M141").Value = -13522.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 3298.8462
$ws.Range("I99").Value = 2488
$ws.Range("J99").Value = 7758.5
$ws.Range("K99").Value = 2488
$ws.Range("L99").Value = 7758.5
$ws.Range("M99").Value = -242
$ws.Range("N99").Value = -12250.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2637.0908
$ws.Range("I46").Value = 879
$ws.Range("K46").Value = 879
$ws.Range("M46").Value = -691

$ws.Range("H61").Value = 6655.3125
$ws.Range("I61").Value = 3927.5
$ws.Range("K61").Value = 3927.5
$ws.Range("M61").Value = -3725.5

$ws.Range("H113").Value = 6655.3125
$ws.Range("I113").Value = 3927.5
$ws.Range("K113").Value = 3927.5
$ws.Range("M113").Value = -1757.5

$ws.Range("H133").Value = 88999
$ws.Range("J133").Value = 88999
$ws.Range("L133").Value = 88999
$ws.Range("N133").Value = -94059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 32140
$ws.Range("J94").Value = 32140
$ws.Range("L94").Value = 32140
$ws.Range("N94").Value = -33942

$ws.Range("H107").Value = 250001170
$ws.Range("I107").Value = 250001300
$ws.Range("K107").Value = 750003900
$ws.Range("M107").Value = -750001980

$ws.Range("H132").Value = 4253
$ws.Range("I132").Value = 3950.6155
$ws.Range("J132").Value = 6218.5
$ws.Range("K132").Value = 11851.8465
$ws.Range("L132").Value = 18655.5
$ws.Range("M132").Value = -9321.8465
$ws.Range("N132").Value = -23715.5

$ws.Range("H133").Value = 44999.668
$ws.Range("J133").Value = 42999.5
$ws.Range("L133").Value = 42999.5
$ws.Range("N133").Value = -53119.5

$ws.Range("H135").Value = 72999
$ws.Range("J135").Value = 72999
$ws.Range("L135").Value = 72999
$ws.Range("N135").Value = -83139

$ws.Range("H136").Value = 18872026
$ws.Range("I136").Value = 25004610
$ws.Range("K136").Value = 75013830
$ws.Range("M136").Value = -75011280
